$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J40").Value = 4750
$ws.Range("I40").Value = 2000
$ws.Range("H40").Value = 3833.3333
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 4750
$ws.Range("N40").Value = -5100
$ws.Range("M40").Value = -1825
$ws.Range("H43").Value = 5112.25
$ws.Range("M43").Value = -1397
$ws.Range("I43").Value = 1466
$ws.Range("K43").Value = 1466
$ws.Range("L74").Value = 5200
$ws.Range("H74").Value = 4999
$ws.Range("J74").Value = 5200
$ws.Range("M74").Value = -3862
$ws.Range("N74").Value = -7072
$ws.Range("I74").Value = 4798
$ws.Range("K74").Value = 4798
$ws.Range("L77").Value = 26000
$ws.Range("H77").Value = 4999
$ws.Range("N77").Value = -35360
$ws.Range("M77").Value = -19310
$ws.Range("K77").Value = 23990
$ws.Range("I77").Value = 4798
$ws.Range("J77").Value = 5200
$ws.Range("J121").Value = 4370.7144
$ws.Range("H121").Value = 4370.7144
$ws.Range("L121").Value = 13112.1432
$ws.Range("N121").Value = -16606.1432
$ws.Range("M135").Value = -11025.5169
$ws.Range("I135").Value = 1506.7241
$ws.Range("L135").Value = 0
$ws.Range("H135").Value = 1506.7241
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13560.5169
$ws.Range("N135").ClearContents()
$ws.Range("M137").Value = -5191.6362
$ws.Range("I137").Value = 2580.5454
$ws.Range("K137").Value = 7741.6362
$ws.Range("H137").Value = 4467.9
$ws.Range("H141").Value = 21388.361
$ws.Range("K141").Value = 71427.288
$ws.Range("M141").Value = -66247.288
$ws.Range("I141").Value = 23809.096

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J2").Value = 3497.5
$ws.Range("M2").Value = -3244.7144
$ws.Range("K2").Value = 3357.7144
$ws.Range("H2").Value = 3388.7778
$ws.Range("L2").Value = 3497.5
$ws.Range("I2").Value = 3357.7144
$ws.Range("N2").Value = -3723.5
$ws.Range("L74").Value = 2133.1667
$ws.Range("H74").Value = 1614.3334
$ws.Range("J74").Value = 2133.1667
$ws.Range("M74").Value = -394.4445000000001
$ws.Range("N74").Value = -3881.1667
$ws.Range("I74").Value = 1268.4445
$ws.Range("K74").Value = 1268.4445
$ws.Range("L77").Value = 10665.8335
$ws.Range("H77").Value = 1614.3334
$ws.Range("N77").Value = -19401.8335
$ws.Range("M77").Value = -1974.2225
$ws.Range("K77").Value = 6342.2225
$ws.Range("I77").Value = 1268.4445
$ws.Range("J77").Value = 2133.1667
$ws.Range("H80").Value = 20000
$ws.Range("N80").Value = -21996
$ws.Range("L80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("N83").Value = -69984
$ws.Range("L83").Value = 60000
$ws.Range("H83").Value = 20000
$ws.Range("N116").Value = -8085.5
$ws.Range("I116").Value = 3357.7144
$ws.Range("M116").Value = -1063.7144
$ws.Range("K116").Value = 3357.7144
$ws.Range("J116").Value = 3497.5
$ws.Range("H116").Value = 3388.7778
$ws.Range("L116").Value = 3497.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 3497.5
$ws.Range("N3").Value = -3725.5
$ws.Range("L3").Value = 3497.5
$ws.Range("H3").Value = 3388.7778
$ws.Range("I3").Value = 3357.7144
$ws.Range("M3").Value = -3243.7144
$ws.Range("K3").Value = 3357.7144
$ws.Range("K20").Value = 15308.875
$ws.Range("M20").Value = -15061.875
$ws.Range("I20").Value = 15308.875
$ws.Range("H20").Value = 16050.353
$ws.Range("L82").Value = 24996.666
$ws.Range("H82").Value = 14908.77
$ws.Range("N82").Value = -25762.666
$ws.Range("J82").Value = 24996.666
$ws.Range("L85").Value = 24996.666
$ws.Range("N85").Value = -27648.666
$ws.Range("J85").Value = 24996.666
$ws.Range("H85").Value = 14908.77
$ws.Range("H86").Value = 2909.6365
$ws.Range("M86").Value = -1306.4285
$ws.Range("I86").Value = 2429.4285
$ws.Range("K86").Value = 2429.4285
$ws.Range("I89").Value = 2429.4285
$ws.Range("M89").Value = -6531.1425
$ws.Range("K89").Value = 12147.1425
$ws.Range("H89").Value = 2909.6365
$ws.Range("N105").Value = -7171.4
$ws.Range("I105").Value = 5973.8887
$ws.Range("L105").Value = 3677.4
$ws.Range("H105").Value = 5474.6523
$ws.Range("K105").Value = 5973.8887
$ws.Range("J105").Value = 3677.4
$ws.Range("M105").Value = -4226.8887

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N58").Value = -2804.25
$ws.Range("I58").Value = 1148.6
$ws.Range("H58").Value = 1583.2609
$ws.Range("M58").Value = -945.5999999999999
$ws.Range("K58").Value = 1148.6
$ws.Range("L58").Value = 2398.25
$ws.Range("J58").Value = 2398.25
$ws.Range("H86").Value = 8124.75
$ws.Range("M86").Value = -5710
$ws.Range("I86").Value = 6833
$ws.Range("K86").Value = 6833
$ws.Range("I89").Value = 6833
$ws.Range("M89").Value = -28549
$ws.Range("K89").Value = 34165
$ws.Range("H89").Value = 8124.75
$ws.Range("H122").Value = 38233.285
$ws.Range("M122").Value = -139373.23
$ws.Range("N122").Value = -20147.5
$ws.Range("J122").Value = 5082.5
$ws.Range("I122").Value = 47274.41
$ws.Range("K122").Value = 141823.23
$ws.Range("L122").Value = 15247.5
$ws.Range("J132").Value = 2724.875
$ws.Range("H132").Value = 1951.1897
$ws.Range("M132").Value = -2952.200000000001
$ws.Range("I132").Value = 1827.4
$ws.Range("N132").Value = -13234.625
$ws.Range("K132").Value = 5482.200000000001
$ws.Range("L132").Value = 8174.625
$ws.Range("H134").Value = 3455.8125
$ws.Range("M134").Value = -7524.2001
$ws.Range("K134").Value = 10059.2001
$ws.Range("I134").Value = 3353.0667
$ws.Range("I136").Value = 1148.6
$ws.Range("J136").Value = 2398.25
$ws.Range("L136").Value = 7194.75
$ws.Range("N136").Value = -12294.75
$ws.Range("K136").Value = 3445.8
$ws.Range("M136").Value = -895.7999999999997
$ws.Range("H136").Value = 1583.2609

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J2").Value = 2552.6667
$ws.Range("M2").Value = -599921.2
$ws.Range("K2").Value = 600034.2
$ws.Range("H2").Value = 63460.812
$ws.Range("L2").Value = 15316.0002
$ws.Range("I2").Value = 100005.7
$ws.Range("N2").Value = -15542.0002
$ws.Range("N5").Value = -3296.6
$ws.Range("I5").Value = 924.7143
$ws.Range("K5").Value = 2774.1429
$ws.Range("H5").Value = 966.1667
$ws.Range("M5").Value = -2662.1429
$ws.Range("L5").Value = 3072.6
$ws.Range("J5").Value = 1024.2
$ws.Range("K7").Value = 204.52941
$ws.Range("H7").Value = 67.210526
$ws.Range("I7").Value = 68.17646999999999
$ws.Range("M7").Value = -92.52940999999998
$ws.Range("J12").Value = 146
$ws.Range("N12").Value = -784
$ws.Range("L12").Value = 438
$ws.Range("I12").Value = 141.625
$ws.Range("K12").Value = 424.875
$ws.Range("M12").Value = -251.875
$ws.Range("H12").Value = 144.25
$ws.Range("L34").Value = 8963.499899999999
$ws.Range("H34").Value = 1697.8182
$ws.Range("J34").Value = 2987.8333
$ws.Range("N34").Value = -9131.499899999999
$ws.Range("M38").Value = 117.2
$ws.Range("N38").ClearContents()
$ws.Range("H38").Value = 76.59999999999999
$ws.Range("J38").Value = 0
$ws.Range("I38").Value = 76.59999999999999
$ws.Range("L38").Value = 0
$ws.Range("K38").Value = 229.8
$ws.Range("H55").Value = 50448.477
$ws.Range("N55").Value = -11544.6
$ws.Range("J55").Value = 3730.2
$ws.Range("L55").Value = 11190.6
$ws.Range("J121").Value = 4079.8333
$ws.Range("H121").Value = 27863382
$ws.Range("L121").Value = 12239.4999
$ws.Range("N121").Value = -14859.4999
$ws.Range("M135").Value = -5787.4287
$ws.Range("I135").Value = 924.7143
$ws.Range("L135").Value = 9217.800000000001
$ws.Range("H135").Value = 966.1667
$ws.Range("J135").Value = 1024.2
$ws.Range("K135").Value = 8322.4287
$ws.Range("N135").Value = -14287.8
$ws.Range("L138").Value = 22506.4995
$ws.Range("H138").Value = 4507.091
$ws.Range("J138").Value = 7502.1665
$ws.Range("N138").Value = -32786.49950000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L107").Value = 3271
$ws.Range("H107").Value = 1862.9333
$ws.Range("J107").Value = 3271
$ws.Range("N107").Value = -7111
$ws.Range("N113").Value = -7815.2856
$ws.Range("L113").Value = 3475.2856
$ws.Range("J113").Value = 3475.2856
$ws.Range("H113").Value = 3033.1333
$ws.Range("H122").Value = 884.5714
$ws.Range("M122").Value = -203.7142000000003
$ws.Range("I122").Value = 884.5714
$ws.Range("K122").Value = 2653.7142

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I16").Value = 4093.1765
$ws.Range("L16").Value = 2000.6666
$ws.Range("M16").Value = -3923.1765
$ws.Range("K16").Value = 4093.1765
$ws.Range("N16").Value = -2340.6666
$ws.Range("J16").Value = 2000.6666
$ws.Range("H16").Value = 3779.3
$ws.Range("K61").Value = 91822.82000000001
$ws.Range("H61").Value = 91822.82000000001
$ws.Range("I61").Value = 91822.82000000001
$ws.Range("M61").Value = -91620.82000000001
$ws.Range("I113").Value = 91822.82000000001
$ws.Range("K113").Value = 91822.82000000001
$ws.Range("M113").Value = -89652.82000000001
$ws.Range("H113").Value = 91822.82000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M107").Value = -5346.333
$ws.Range("I107").Value = 2422.111
$ws.Range("H107").Value = 38463816
$ws.Range("K107").Value = 7266.333
$ws.Range("H132").Value = 975775
$ws.Range("M132").Value = -12610.2281
$ws.Range("I132").Value = 5046.7427
$ws.Range("K132").Value = 15140.2281
